$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the Start time / Completion time of the first three sample flight
# rows forward by 9 days (date only - time of day is unchanged).
$ws.Range("B2").Value = 43475.647037037037
$ws.Range("C2").Value = 43475.647893518515
$ws.Range("B3").Value = 43475.647974537038
$ws.Range("C3").Value = 43475.648449074077
$ws.Range("B4").Value = 43475.648472222223
$ws.Range("C4").Value = 43475.648877314816

# Update the active cell/selection saved with the sheet view.
$ws.Range("C14").Select()
